# Database Grading Criteria - "updated comments and criteria"
#
# Semantic changes (decoded from the OOXML diff):
#   - Column C ("*" marker / comment column) toggled on a handful of rows:
#       C5  (Update)                         : blank    -> "*"
#       C8  (Correctly formatted SQL File)    : blank    -> "*"
#       C10 (ALTER TABLE)                     : "arguable" -> "*"
#       C15 (View)                            : "*"      -> blank
#       C16 (Sequence)                        : blank    -> "*"
#       C25 (WHERE)                           : blank    -> "*"
#       C26 (AND/OR)                          : blank    -> "*"
#   - The now-unused shared string "arguable" disappears on save once no
#     cell references it any more (handled automatically by the engine).
#   - The saved cursor/selection moves from C17 to C29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "*"
$ws.Range("C8").Value = "*"
$ws.Range("C10").Value = "*"
$ws.Range("C15").ClearContents()
$ws.Range("C16").Value = "*"
$ws.Range("C25").Value = "*"
$ws.Range("C26").Value = "*"

$ws.Range("C29").Select()
